$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-unused rows 17:21 (table shrinks from 20 to 15 data rows)
$ws.Rows("17:21").Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl25"
$ws.Cells.Item(2, 3).Value = "Ackr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 6.91741
$ws.Cells.Item(2, 8).Value = 20.75223
$ws.Cells.Item(2, 9).Value = 0.2334435312127427
$ws.Cells.Item(2, 10).Value = 0.2334435312127427
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.377371
$ws.Cells.Item(2, 14).Value = 1.132113
$ws.Cells.Item(2, 15).Value = 0.4698794580655765
$ws.Cells.Item(2, 16).Value = 0.4698794580655764
$ws.Cells.Item(2, 17).Value = 2.61042992911
$ws.Cells.Item(2, 18).Value = 23.49386936199
$ws.Cells.Item(2, 19).Value = 0.109690319935158
$ws.Cells.Item(2, 20).Value = 0.109690319935158

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl25"
$ws.Cells.Item(3, 3).Value = "Ackr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 6.91741
$ws.Cells.Item(3, 8).Value = 20.75223
$ws.Cells.Item(3, 9).Value = 0.2334435312127427
$ws.Cells.Item(3, 10).Value = 0.2334435312127427
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.3560133333333333
$ws.Cells.Item(3, 14).Value = 1.06804
$ws.Cells.Item(3, 15).Value = 0.443286188209444
$ws.Cells.Item(3, 16).Value = 0.443286188209444
$ws.Cells.Item(3, 17).Value = 2.462690192133334
$ws.Cells.Item(3, 18).Value = 22.1642117292
$ws.Cells.Item(3, 19).Value = 0.1034822931134491
$ws.Cells.Item(3, 20).Value = 0.1034822931134491

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ccl25"
$ws.Cells.Item(4, 3).Value = "Ackr4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.91741
$ws.Cells.Item(4, 8).Value = 20.75223
$ws.Cells.Item(4, 9).Value = 0.2334435312127427
$ws.Cells.Item(4, 10).Value = 0.2334435312127427
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.06973866666666667
$ws.Cells.Item(4, 14).Value = 0.209216
$ws.Cells.Item(4, 15).Value = 0.08683435372497944
$ws.Cells.Item(4, 16).Value = 0.08683435372497944
$ws.Cells.Item(4, 17).Value = 0.4824109501866667
$ws.Cells.Item(4, 18).Value = 4.34169855168
$ws.Cells.Item(4, 19).Value = 0.02027091816413558
$ws.Cells.Item(4, 20).Value = 0.02027091816413558

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl25"
$ws.Cells.Item(5, 3).Value = "Ackr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 7.803964999999999
$ws.Cells.Item(5, 8).Value = 23.411895
$ws.Cells.Item(5, 9).Value = 0.2633623201546029
$ws.Cells.Item(5, 10).Value = 0.2633623201546028
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.377371
$ws.Cells.Item(5, 14).Value = 1.132113
$ws.Cells.Item(5, 15).Value = 0.4698794580655765
$ws.Cells.Item(5, 16).Value = 0.4698794580655764
$ws.Cells.Item(5, 17).Value = 2.944990076014999
$ws.Cells.Item(5, 18).Value = 26.50491068413499
$ws.Cells.Item(5, 19).Value = 0.1237485442691376
$ws.Cells.Item(5, 20).Value = 0.1237485442691376

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ccl25"
$ws.Cells.Item(6, 3).Value = "Ackr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7.803964999999999
$ws.Cells.Item(6, 8).Value = 23.411895
$ws.Cells.Item(6, 9).Value = 0.2633623201546029
$ws.Cells.Item(6, 10).Value = 0.2633623201546028
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.3560133333333333
$ws.Cells.Item(6, 14).Value = 1.06804
$ws.Cells.Item(6, 15).Value = 0.443286188209444
$ws.Cells.Item(6, 16).Value = 0.443286188209444
$ws.Cells.Item(6, 17).Value = 2.778315592866667
$ws.Cells.Item(6, 18).Value = 25.0048403358
$ws.Cells.Item(6, 19).Value = 0.1167448790193291
$ws.Cells.Item(6, 20).Value = 0.1167448790193291

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ccl25"
$ws.Cells.Item(7, 3).Value = "Ackr4"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 7.803964999999999
$ws.Cells.Item(7, 8).Value = 23.411895
$ws.Cells.Item(7, 9).Value = 0.2633623201546029
$ws.Cells.Item(7, 10).Value = 0.2633623201546028
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.06973866666666667
$ws.Cells.Item(7, 14).Value = 0.209216
$ws.Cells.Item(7, 15).Value = 0.08683435372497944
$ws.Cells.Item(7, 16).Value = 0.08683435372497944
$ws.Cells.Item(7, 17).Value = 0.5442381138133333
$ws.Cells.Item(7, 18).Value = 4.898143024319999
$ws.Cells.Item(7, 19).Value = 0.02286889686613607
$ws.Cells.Item(7, 20).Value = 0.02286889686613606

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Ccl25"
$ws.Cells.Item(8, 3).Value = "Ackr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 6.430676666666667
$ws.Cells.Item(8, 8).Value = 19.29203
$ws.Cells.Item(8, 9).Value = 0.2170176220802376
$ws.Cells.Item(8, 10).Value = 0.2170176220802376
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.377371
$ws.Cells.Item(8, 14).Value = 1.132113
$ws.Cells.Item(8, 15).Value = 0.4698794580655765
$ws.Cells.Item(8, 16).Value = 0.4698794580655764
$ws.Cells.Item(8, 17).Value = 2.426750884376666
$ws.Cells.Item(8, 18).Value = 21.84075795939
$ws.Cells.Item(8, 19).Value = 0.1019721226537421
$ws.Cells.Item(8, 20).Value = 0.1019721226537421

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Ccl25"
$ws.Cells.Item(9, 3).Value = "Ackr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 6.430676666666667
$ws.Cells.Item(9, 8).Value = 19.29203
$ws.Cells.Item(9, 9).Value = 0.2170176220802376
$ws.Cells.Item(9, 10).Value = 0.2170176220802376
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.3560133333333333
$ws.Cells.Item(9, 14).Value = 1.06804
$ws.Cells.Item(9, 15).Value = 0.443286188209444
$ws.Cells.Item(9, 16).Value = 0.443286188209444
$ws.Cells.Item(9, 17).Value = 2.289406635688889
$ws.Cells.Item(9, 18).Value = 20.6046597212
$ws.Cells.Item(9, 19).Value = 0.09620091446622618
$ws.Cells.Item(9, 20).Value = 0.09620091446622618

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Ccl25"
$ws.Cells.Item(10, 3).Value = "Ackr4"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 6.430676666666667
$ws.Cells.Item(10, 8).Value = 19.29203
$ws.Cells.Item(10, 9).Value = 0.2170176220802376
$ws.Cells.Item(10, 10).Value = 0.2170176220802376
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.06973866666666667
$ws.Cells.Item(10, 14).Value = 0.209216
$ws.Cells.Item(10, 15).Value = 0.08683435372497944
$ws.Cells.Item(10, 16).Value = 0.08683435372497944
$ws.Cells.Item(10, 17).Value = 0.4484668164977778
$ws.Cells.Item(10, 18).Value = 4.036201348480001
$ws.Cells.Item(10, 19).Value = 0.01884458496026926
$ws.Cells.Item(10, 20).Value = 0.01884458496026926

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ccl25"
$ws.Cells.Item(11, 3).Value = "Ackr4"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.144241333333333
$ws.Cells.Item(11, 8).Value = 12.432724
$ws.Cells.Item(11, 9).Value = 0.1398567283204463
$ws.Cells.Item(11, 10).Value = 0.1398567283204463
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.377371
$ws.Cells.Item(11, 14).Value = 1.132113
$ws.Cells.Item(11, 15).Value = 0.4698794580655765
$ws.Cells.Item(11, 16).Value = 0.4698794580655764
$ws.Cells.Item(11, 17).Value = 1.563916496201333
$ws.Cells.Item(11, 18).Value = 14.075248465812
$ws.Cells.Item(11, 19).Value = 0.06571580371003587
$ws.Cells.Item(11, 20).Value = 0.06571580371003585

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ccl25"
$ws.Cells.Item(12, 3).Value = "Ackr4"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.144241333333333
$ws.Cells.Item(12, 8).Value = 12.432724
$ws.Cells.Item(12, 9).Value = 0.1398567283204463
$ws.Cells.Item(12, 10).Value = 0.1398567283204463
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.3560133333333333
$ws.Cells.Item(12, 14).Value = 1.06804
$ws.Cells.Item(12, 15).Value = 0.443286188209444
$ws.Cells.Item(12, 16).Value = 0.443286188209444
$ws.Cells.Item(12, 17).Value = 1.475405171217778
$ws.Cells.Item(12, 18).Value = 13.27864654096
$ws.Cells.Item(12, 19).Value = 0.06199655599261444
$ws.Cells.Item(12, 20).Value = 0.06199655599261444

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ccl25"
$ws.Cells.Item(13, 3).Value = "Ackr4"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.144241333333333
$ws.Cells.Item(13, 8).Value = 12.432724
$ws.Cells.Item(13, 9).Value = 0.1398567283204463
$ws.Cells.Item(13, 10).Value = 0.1398567283204463
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.06973866666666667
$ws.Cells.Item(13, 14).Value = 0.209216
$ws.Cells.Item(13, 15).Value = 0.08683435372497944
$ws.Cells.Item(13, 16).Value = 0.08683435372497944
$ws.Cells.Item(13, 17).Value = 0.2890138649315556
$ws.Cells.Item(13, 18).Value = 2.601124784384
$ws.Cells.Item(13, 19).Value = 0.01214436861779598
$ws.Cells.Item(13, 20).Value = 0.01214436861779598

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Ccl25"
$ws.Cells.Item(14, 3).Value = "Ackr4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 4.335755333333334
$ws.Cells.Item(14, 8).Value = 13.007266
$ws.Cells.Item(14, 9).Value = 0.1463197982319706
$ws.Cells.Item(14, 10).Value = 0.1463197982319706
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.377371
$ws.Cells.Item(14, 14).Value = 1.132113
$ws.Cells.Item(14, 15).Value = 0.4698794580655765
$ws.Cells.Item(14, 16).Value = 0.4698794580655764
$ws.Cells.Item(14, 17).Value = 1.636188325895333
$ws.Cells.Item(14, 18).Value = 14.725694933058
$ws.Cells.Item(14, 19).Value = 0.06875266749750285
$ws.Cells.Item(14, 20).Value = 0.06875266749750282

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Ccl25"
$ws.Cells.Item(15, 3).Value = "Ackr4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 4.335755333333334
$ws.Cells.Item(15, 8).Value = 13.007266
$ws.Cells.Item(15, 9).Value = 0.1463197982319706
$ws.Cells.Item(15, 10).Value = 0.1463197982319706
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 0.3333333333333333
$ws.Cells.Item(15, 13).Value = 0.3560133333333333
$ws.Cells.Item(15, 14).Value = 1.06804
$ws.Cells.Item(15, 15).Value = 0.443286188209444
$ws.Cells.Item(15, 16).Value = 0.443286188209444
$ws.Cells.Item(15, 17).Value = 1.543586708737778
$ws.Cells.Item(15, 18).Value = 13.89228037864
$ws.Cells.Item(15, 19).Value = 0.0648615456178252
$ws.Cells.Item(15, 20).Value = 0.06486154561782519

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Ccl25"
$ws.Cells.Item(16, 3).Value = "Ackr4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 4.335755333333334
$ws.Cells.Item(16, 8).Value = 13.007266
$ws.Cells.Item(16, 9).Value = 0.1463197982319706
$ws.Cells.Item(16, 10).Value = 0.1463197982319706
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.06973866666666667
$ws.Cells.Item(16, 14).Value = 0.209216
$ws.Cells.Item(16, 15).Value = 0.08683435372497944
$ws.Cells.Item(16, 16).Value = 0.08683435372497944
$ws.Cells.Item(16, 17).Value = 0.3023697959395557
$ws.Cells.Item(16, 18).Value = 2.721328163456
$ws.Cells.Item(16, 19).Value = 0.01270558511664256
$ws.Cells.Item(16, 20).Value = 0.01270558511664256
